$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build a 2D array (rows x cols) for the new data block A2:I25
$nRows = 24
$nCols = 9
$arr = New-Object 'object[,]' $nRows,$nCols

$arr[0,0] = "model_8_0_0"
$arr[0,1] = 0.00360789897647229
$arr[0,2] = -0.02227312218359745
$arr[0,3] = -0.6049943525983172
$arr[0,4] = -0.04654398098850288
$arr[0,5] = 1.102712392807007
$arr[0,6] = 2.424077033996582
$arr[0,7] = 0.2013064026832581
$arr[0,8] = 1.378067493438721

$arr[1,0] = "model_8_0_23"
$arr[1,1] = 0.2822280652317948
$arr[1,2] = -0.5309069578280707
$arr[1,3] = -1.094049999689128
$arr[1,4] = -0.5533693074731696
$arr[1,5] = 0.7943620681762695
$arr[1,6] = 3.630181312561035
$arr[1,7] = 0.2626461982727051
$arr[1,8] = 2.045444250106812

$arr[2,0] = "model_8_0_22"
$arr[2,1] = 0.2907712234455657
$arr[2,2] = -0.5081051577113975
$arr[2,3] = -1.048414813234539
$arr[2,4] = -0.529595864114093
$arr[2,5] = 0.7849072813987732
$arr[2,6] = 3.576112031936646
$arr[2,7] = 0.2569224238395691
$arr[2,8] = 2.014140129089355

$arr[3,0] = "model_8_0_21"
$arr[3,1] = 0.294137774843938
$arr[3,2] = -0.4981981367846011
$arr[3,3] = -1.048031997947845
$arr[3,4] = -0.5201242567299311
$arr[3,5] = 0.7811815142631531
$arr[3,6] = 3.552619934082031
$arr[3,7] = 0.2568743824958801
$arr[3,8] = 2.001668214797974

$arr[4,0] = "model_8_0_1"
$arr[4,1] = 0.3004109495316379
$arr[4,2] = 0.1768815845368426
$arr[4,3] = -0.794122486585674
$arr[4,4] = 0.1348460141673719
$arr[4,5] = 0.7742389440536499
$arr[4,6] = 1.951829075813293
$arr[4,7] = 0.2250277996063232
$arr[4,8] = 1.139216780662537

$arr[5,0] = "model_8_0_20"
$arr[5,1] = 0.4574488464323029
$arr[5,2] = -0.04812893193633183
$arr[5,3] = -0.4707736330875349
$arr[5,4] = -0.06518828267963328
$arr[5,5] = 0.6004442572593689
$arr[5,6] = 2.485388040542603
$arr[5,7] = 0.1844717711210251
$arr[5,8] = 1.402617812156677

$arr[6,0] = "model_8_0_2"
$arr[6,1] = 0.5992330177109119
$arr[6,2] = 0.4746546718295044
$arr[6,3] = -0.7919599252067782
$arr[6,4] = 0.4188300163124785
$arr[6,5] = 0.4435309767723083
$arr[6,6] = 1.245731234550476
$arr[6,7] = 0.2247565537691116
$arr[6,8] = 0.7652725577354431

$arr[7,0] = "model_8_0_19"
$arr[7,1] = 0.6139026109794967
$arr[7,2] = 0.4037457942014204
$arr[7,3] = -0.3774078374511685
$arr[7,4] = 0.3698023380339559
$arr[7,5] = 0.4272960722446442
$arr[7,6] = 1.413874864578247
$arr[7,7] = 0.1727613806724548
$arr[7,8] = 0.8298312425613403

$arr[8,0] = "model_8_0_3"
$arr[8,1] = 0.6315510331374967
$arr[8,2] = 0.4989189677322526
$arr[8,3] = -0.3593531004340638
$arr[8,4] = 0.4613541364988543
$arr[8,5] = 0.4077644348144531
$arr[8,6] = 1.188194155693054
$arr[8,7] = 0.1704968512058258
$arr[8,8] = 0.7092776894569397

$arr[9,0] = "model_8_0_4"
$arr[9,1] = 0.6446089610308425
$arr[9,2] = 0.5246546061665023
$arr[9,3] = -0.6400777723154734
$arr[9,4] = 0.4733062681158334
$arr[9,5] = 0.3933131396770477
$arr[9,6] = 1.127168297767639
$arr[9,7] = 0.2057067453861237
$arr[9,8] = 0.6935393810272217

$arr[10,0] = "model_8_0_5"
$arr[10,1] = 0.6478333346487901
$arr[10,2] = 0.5266482035835585
$arr[10,3] = -0.5931010773344256
$arr[10,4] = 0.4773117046553086
$arr[10,5] = 0.3897447288036346
$arr[10,6] = 1.122440934181213
$arr[10,7] = 0.1998146921396255
$arr[10,8] = 0.6882650852203369

$arr[11,0] = "model_8_0_18"
$arr[11,1] = 0.6488529618231931
$arr[11,2] = 0.5036245067058069
$arr[11,3] = -0.3243181726998603
$arr[11,4] = 0.4674128882614131
$arr[11,5] = 0.3886162936687469
$arr[11,6] = 1.177036285400391
$arr[11,7] = 0.1661026030778885
$arr[11,8] = 0.7012996077537537

$arr[12,0] = "model_8_0_6"
$arr[12,1] = 0.6488895221339841
$arr[12,2] = 0.5289071697092835
$arr[12,3] = -0.6048718333543961
$arr[12,4] = 0.4789394979570819
$arr[12,5] = 0.3885758221149445
$arr[12,6] = 1.117084383964539
$arr[12,7] = 0.2012910395860672
$arr[12,8] = 0.6861215829849243

$arr[13,0] = "model_8_0_7"
$arr[13,1] = 0.6496967673965036
$arr[13,2] = 0.5304208449694559
$arr[13,3] = -0.6033388722990221
$arr[13,4] = 0.480449939744062
$arr[13,5] = 0.3876824378967285
$arr[13,6] = 1.113495111465454
$arr[13,7] = 0.201098769903183
$arr[13,8] = 0.6841326951980591

$arr[14,0] = "model_8_0_8"
$arr[14,1] = 0.6501421426076986
$arr[14,2] = 0.5311499282923202
$arr[14,3] = -0.6040755523908299
$arr[14,4] = 0.481112427915852
$arr[14,5] = 0.3871895968914032
$arr[14,6] = 1.111766219139099
$arr[14,7] = 0.2011911571025848
$arr[14,8] = 0.683260440826416

$arr[15,0] = "model_8_0_9"
$arr[15,1] = 0.6506746957982494
$arr[15,2] = 0.5323514239386353
$arr[15,3] = -0.6153923569655904
$arr[15,4] = 0.4817507572294563
$arr[15,5] = 0.3866001665592194
$arr[15,6] = 1.108917117118835
$arr[15,7] = 0.202610582113266
$arr[15,8] = 0.6824198365211487

$arr[16,0] = "model_8_0_10"
$arr[16,1] = 0.6514093954526898
$arr[16,2] = 0.5328093295290721
$arr[16,3] = -0.591919651979951
$arr[16,4] = 0.4832393984538995
$arr[16,5] = 0.3857870697975159
$arr[16,6] = 1.107831478118896
$arr[16,7] = 0.1996665000915527
$arr[16,8] = 0.6804595589637756

$arr[17,0] = "model_8_0_17"
$arr[17,1] = 0.6533029526046494
$arr[17,2] = 0.518142247252275
$arr[17,3] = -0.3487339646701004
$arr[17,4] = 0.4801577791557413
$arr[17,5] = 0.3836914598941803
$arr[17,6] = 1.142611026763916
$arr[17,7] = 0.1691649556159973
$arr[17,8] = 0.6845174431800842

$arr[18,0] = "model_8_0_14"
$arr[18,1] = 0.6537689393600604
$arr[18,2] = 0.5288757486652287
$arr[18,3] = -0.5342207618141623
$arr[18,4] = 0.4820750780991628
$arr[18,5] = 0.383175790309906
$arr[18,6] = 1.117158889770508
$arr[18,7] = 0.1924296319484711
$arr[18,8] = 0.681992769241333

$arr[19,0] = "model_8_0_16"
$arr[19,1] = 0.6537974102497708
$arr[19,2] = 0.5228622695504979
$arr[19,3] = -0.4139467010893219
$arr[19,4] = 0.4817339519205807
$arr[19,5] = 0.3831442296504974
$arr[19,6] = 1.131418466567993
$arr[19,7] = 0.1773442625999451
$arr[19,8] = 0.6824420094490051

$arr[20,0] = "model_8_0_15"
$arr[20,1] = 0.6543678848162855
$arr[20,2] = 0.5272281078528994
$arr[20,3] = -0.4622403046653702
$arr[20,4] = 0.483730834077388
$arr[20,5] = 0.3825128972530365
$arr[20,6] = 1.121065855026245
$arr[20,7] = 0.1834014803171158
$arr[20,8] = 0.6798125505447388

$arr[21,0] = "model_8_0_13"
$arr[21,1] = 0.6549559739567312
$arr[21,2] = 0.5334523394599244
$arr[21,3] = -0.5480448412247481
$arr[21,4] = 0.4858191511786049
$arr[21,5] = 0.3818620145320892
$arr[21,6] = 1.106306552886963
$arr[21,7] = 0.194163516163826
$arr[21,8] = 0.6770626306533813

$arr[22,0] = "model_8_0_11"
$arr[22,1] = 0.6574526976448981
$arr[22,2] = 0.5346935729252333
$arr[22,3] = -0.3270949505228409
$arr[22,4] = 0.4969063338771953
$arr[22,5] = 0.3790989518165588
$arr[22,6] = 1.103363275527954
$arr[22,7] = 0.1664508730173111
$arr[22,8] = 0.6624633073806763

$arr[23,0] = "model_8_0_12"
$arr[23,1] = 0.6589533625214915
$arr[23,2] = 0.5296800645862065
$arr[23,3] = -0.1852668316409705
$arr[23,4] = 0.4984841175765918
$arr[23,5] = 0.3774381279945374
$arr[23,6] = 1.115251541137695
$arr[23,7] = 0.1486620754003525
$arr[23,8] = 0.6603856682777405

$ws.Range("A2:I25").Value = $arr

# Row 26 no longer exists in the new data; clear leftover old row
$ws.Range("A26:I26").Clear()

$ws.Range("A1").Select()
